$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.068.37'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '1.830.76'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.62'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.008'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4619'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3702'
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07353'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8732'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07974'
$ws.Range('E11').Value = '  +3.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.83'
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.353'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.558'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.91'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.725.14'
$ws.Range('E16').Value = '  -9.15%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008876'
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.70'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').Value = '26.908.05'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.59'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '1.992.57'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.67'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.848'
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.58'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.077'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.087'
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.47'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08885'
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.979'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7328'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.137'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.456'
$ws.Range('E36').Value = '  -3.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.072'
$ws.Range('E37').Value = '  -1.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01947'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05252'
$ws.Range('E39').Value = '  -0.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.947'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.143'
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5180'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1632'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8606'
$ws.Range('E44').Value = '  -14.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.227'
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4833'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.28'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.009'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.23'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06234'
$ws.Range('E51').Value = '  -0.78%  '
